# Update column F (dSF) values for the rodón_carlos.xlsx data sheet.
# The commit message indicates data was re-pulled and the "dSF" (delta-score-final)
# column was recalculated for each game row. This script pushes the new values
# into column F for the affected rows, leaving all other data untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -3
    3  = -2
    4  = 3
    5  = -1
    7  = 6
    8  = 2
    9  = -1
    10 = -1
    11 = -3
    12 = 2
    13 = -4
    14 = 7
    15 = 5
    17 = 4
    18 = -1
    19 = -3
    20 = -3
    21 = -6
    22 = -6
    23 = -3
    24 = 3
    25 = 6
    26 = 6
    27 = 7
    28 = 5
    29 = 4
    30 = 7
    31 = -3
    32 = 10
    34 = -2
    35 = 3
    37 = -1
    38 = -3
    39 = 3
    41 = -1
    42 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
